# Arvados AWS install checklist: add a "ClusterID" row near the top of
# the worksheet, just above "Machines".
#
# The original sheet carried a stray "last row" marker at row 1048576
# (an artifact of the prior export) - remove that first so that
# inserting new rows doesn't push it out of the valid row range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the stray trailing sentinel row before shifting everything down.
$ws.Rows(1048576).Delete()

# Insert two new rows at the top of the checklist (row 3/4), pushing
# "Machines" and everything below it down by two rows, matching the
# blank-row spacing used elsewhere in the sheet.
$ws.Rows("3:4").Insert()

# Populate the newly inserted heading row.
$ws.Range("A3").Value = "ClusterID"

# Nudge the bottom-right corner cell so the sheet's used range / saved
# dimension extends down through the final (blank) row 47, matching the
# trimmed sheet extent after the edit.
$ws.Range("E47").NumberFormat = "General"

# Restore the expected active selection/cursor position.
$ws.Range("A4").Select() | Out-Null
